$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 7, pushing existing rows 7-9 down to 8-10.
$ws.Rows.Item(7).Insert()

# Populate the new row 7 with the new record (Damasco, variedad "Modesto").
$ws.Cells.Item(7, 1).Value = 11
$ws.Cells.Item(7, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(7, 3).Value = "Bíobío"
$ws.Cells.Item(7, 4).Value = 44579
$ws.Cells.Item(7, 4).Style = $ws.Cells.Item(8, 4).Style
$ws.Cells.Item(7, 4).NumberFormat = $ws.Cells.Item(8, 4).NumberFormat
$ws.Cells.Item(7, 5).Value = 8
$ws.Cells.Item(7, 6).Value = "Fruta"
$ws.Cells.Item(7, 7).Value = 100103
$ws.Cells.Item(7, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(7, 9).Value = 100103003
$ws.Cells.Item(7, 10).Value = "Damasco"
$ws.Cells.Item(7, 11).Value = "Modesto"
$ws.Cells.Item(7, 12).Value = "Primera"
$ws.Cells.Item(7, 13).Value = 180
$ws.Cells.Item(7, 14).Value = 13000
$ws.Cells.Item(7, 15).Value = 14000
$ws.Cells.Item(7, 16).Value = 13444
$ws.Cells.Item(7, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(7, 18).Value = "Región Metropolitana"
$ws.Cells.Item(7, 19).Value = 747
$ws.Cells.Item(7, 20).Value = 18
